# Update the gene symbol list for the L5/6 Neurons row to include the
# newly added marker genes C1QL3 and SCN4B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "THEMIS, FEZF2, C1QL3, SCN4B"

# The longer text no longer fits the previously "best fit" column width,
# so widen column C (geneSymbolmore1) to accommodate it.
$ws.Columns.Item(3).ColumnWidth = 25.33

# Move/update the active selection on the sheet.
$ws.Range("F7").Select()
